$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "330.87"
Set-TextValue $ws "E2" "-0.16%"
Set-TextValue $ws "D3" "43.30"
Set-TextValue $ws "E3" "3.67%"
Set-TextValue $ws "D4" "5.603"
Set-TextValue $ws "E4" "-1.87%"
Set-TextValue $ws "D5" "0.08195"
Set-TextValue $ws "E5" "-2.14%"
Set-TextValue $ws "D6" "8.749"
Set-TextValue $ws "E6" "-0.54%"
Set-TextValue $ws "D7" "4.389"
Set-TextValue $ws "E7" "-3.36%"
Set-TextValue $ws "D8" "1.912"
Set-TextValue $ws "E8" "-5.56%"
Set-TextValue $ws "D10" "0.9439"
Set-TextValue $ws "E10" "1.96%"
Set-TextValue $ws "D11" "0.1198"
Set-TextValue $ws "E11" "-7.63%"
Set-TextValue $ws "D12" "0.1923"
Set-TextValue $ws "E12" "-2.41%"
Set-TextValue $ws "D13" "0.09765"
Set-TextValue $ws "E13" "3.59%"
Set-TextValue $ws "D14" "0.04340"
Set-TextValue $ws "E14" "11.25%"
Set-TextValue $ws "D15" "0.1070"
Set-TextValue $ws "E15" "0.76%"
Set-TextValue $ws "D16" "0.001283"
Set-TextValue $ws "E16" "-1.84%"
Set-TextValue $ws "D17" "0.006014"
Set-TextValue $ws "E17" "-2.22%"
Set-TextValue $ws "D18" "3.512"
Set-TextValue $ws "E18" "2.05%"
Set-TextValue $ws "D19" "0.3537"
Set-TextValue $ws "D20" "8.709"
Set-TextValue $ws "E20" "5.47%"
Set-TextValue $ws "D21" "0.1370"
Set-TextValue $ws "E21" "-0.21%"
Set-TextValue $ws "D23" "0.04392"
Set-TextValue $ws "E23" "-0.55%"
Set-TextValue $ws "D24" "0.001239"
Set-TextValue $ws "E24" "-1.25%"
Set-TextValue $ws "E25" "-1.28%"
Set-TextValue $ws "D26" "0.0001235"
Set-TextValue $ws "D27" "0.0004005"
Set-TextValue $ws "E27" "31.52%"
Set-TextValue $ws "D39" "0.02769"
Set-TextValue $ws "E39" "-1.26%"
Set-TextValue $ws "D40" "0.05732"
Set-TextValue $ws "E40" "3.33%"
Set-TextValue $ws "D41" "0.007922"
Set-TextValue $ws "E41" "1.74%"
Set-TextValue $ws "D42" "0.009770"
Set-TextValue $ws "E42" "5.95%"
Set-TextValue $ws "D43" "0.1417"
Set-TextValue $ws "E43" "-1.23%"
Set-TextValue $ws "E44" "-1.62%"
Set-TextValue $ws "D45" "0.009673"
Set-TextValue $ws "E45" "-12.75%"
Set-TextValue $ws "D46" "0.00007354"
Set-TextValue $ws "E46" "4.86%"
Set-TextValue $ws "D47" "0.00000000753"
Set-TextValue $ws "E47" "0.40%"
Set-TextValue $ws "D48" "0.003446"
Set-TextValue $ws "D49" "0.002279"
Set-TextValue $ws "E49" "-0.04%"
Set-TextValue $ws "D50" "0.00002109"
Set-TextValue $ws "E50" "0.40%"
Set-TextValue $ws "D51" "0.0002008"
Set-TextValue $ws "E51" "0.40%"
